$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.753.17'
$ws.Range('E2').Value = '  -1.50%  '
$ws.Range('D3').Value = '2.541.66'
$ws.Range('E3').Value = '  +2.40%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '566.13'
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.70'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -2.08%  '
$ws.Range('D9').Value = '2.541.93'
$ws.Range('E9').Value = '  +2.44%  '
$ws.Range('E10').Value = '  -2.21%  '
$ws.Range('E11').Value = '  -2.67%  '
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.352'
$ws.Range('E13').Value = '  -1.07%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.07'
$ws.Range('E14').Value = '  +2.05%  '
$ws.Range('D15').Value = '2.990.18'
$ws.Range('E15').Value = '  +2.19%  '
$ws.Range('D16').Value = '62.686.35'
$ws.Range('E17').Value = '  -1.98%  '
$ws.Range('D18').Value = '2.538.67'
$ws.Range('E18').Value = '  +2.23%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.43'
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '333.38'
$ws.Range('E20').Value = '  -2.79%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.27'
$ws.Range('E21').Value = '  -1.44%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.73'
$ws.Range('E22').Value = '  -2.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '64.50'
$ws.Range('E24').Value = '  -2.18%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.169'
$ws.Range('E25').Value = '  -3.98%  '
$ws.Range('E26').Value = '  +3.92%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('B28').Value = 'SuiNetwork'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.48'
$ws.Range('E28').Value = '  +10.58%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.30'
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.23'
$ws.Range('E30').Value = '  +5.16%  '
$ws.Range('D31').Value = '0.0₃0805'
$ws.Range('E31').Value = '  -2.99%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.85'
$ws.Range('E32').Value = '  -1.40%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '176.56'
$ws.Range('E33').Value = '  -0.52%  '
$ws.Range('E34').Value = '  +3.48%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '405.47'
$ws.Range('E35').Value = '  +8.72%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.397'
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.85'
$ws.Range('E37').Value = '  -0.96%  '
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.33'
$ws.Range('E39').Value = '  -3.51%  '
$ws.Range('E40').Value = '  +1.25%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '38.96'
$ws.Range('E42').Value = '  -3.70%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '152.64'
$ws.Range('E43').Value = '  +0.73%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.73'
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '20.59'
$ws.Range('E45').Value = '  -1.63%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.602'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0954'
$ws.Range('E47').Value = '  -1.50%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0516'
$ws.Range('E48').Value = '  -1.93%  '
$ws.Range('E49').Value = '  +3.08%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '18.16'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.76'
$ws.Range('E51').Value = '  -2.33%  '
